$wb = $excel.ActiveWorkbook

# --- Sheet: Estadisticos 2P (row 2 and row 6 stats update) ---
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Cells.Item(2, 4).Value = 15
$ws2.Cells.Item(2, 5).Value = 7
$ws2.Cells.Item(2, 6).Value = 26
$ws2.Cells.Item(2, 7).Value = 63.41
$ws2.Cells.Item(2, 8).Value = 7.8

$ws2.Cells.Item(6, 4).Value = 18
$ws2.Cells.Item(6, 5).Value = 3
$ws2.Cells.Item(6, 6).Value = 18
$ws2.Cells.Item(6, 7).Value = 50
$ws2.Cells.Item(6, 8).Value = 8.300000000000001

# --- Sheet: Estadisticos Final (row 2 and row 6 promedio update) ---
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Cells.Item(2, 8).Value = 8
$ws3.Cells.Item(6, 8).Value = 8.4

# --- Sheet: Rescatables (append student rows) ---
$ws4 = $wb.Worksheets.Item("Rescatables")

$rescatables = @(
    ,@(2, 20330051920374, "BERISTAIN", "APALE", "JOSE ISAIAS", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "2ALCM", 2)
    ,@(3, 20330051920223, "CASTELLANOS", "TEQUIHUATLE", "JENNIFER", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "2ALCM", 2)
    ,@(4, 20330051920233, "LASTRE", "PACHECO", "ATENEA", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "2ALCM", 2)
    ,@(5, 20330051920238, "MARTINEZ", "CARRERA", "CAROLINA", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "2ALCM", 2)
    ,@(6, 20330051920252, "SANCHEZ", "PEREZ", "ARLET", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "2ALCM", 2)
    ,@(7, 20330051920259, "XOTLANIHUA", "XOTLANIHUA", "JESUS URIEL", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "2ALCM", 2)
    ,@(8, 20330051920341, "CRUZ", "PIMENTEL", "ISYSS MONSERRATH", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "2APM", 2)
    ,@(9, 20330051920343, "GERARDO", "CASTRO", "EDUARDO", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "2APM", 2)
    ,@(10, 20330051920344, "HERNANDEZ", "AGUIRRE", "MARIAM GUADALUPE", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "2APM", 2)
    ,@(11, 20330051920345, "HERNANDEZ", "IXTLA", "ELIEL", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "2APM", 2)
    ,@(12, 20330051920348, "MIRANDA", "ESTRELLA", "JESUS", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "2APM", 2)
    ,@(13, 20330051920349, "NAMIGTLE", "MOLOHUA", "ANGEL ALDAHIR", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "2APM", 2)
    ,@(14, 20330051920350, "ROMAN", "HERNANDEZ", "ESTEFANI", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "2APM", 2)
    ,@(15, 20330051920352, "SANCHEZ", "RODRIGUEZ", "EMILIO", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "2APM", 2)
    ,@(16, 20330051920353, "TORRES", "PEREZ", "INGRID", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "2APM", 2)
    ,@(17, 20330051920354, "TRUJILLO", "HERRERA", "JAHEL", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "2APM", 2)
    ,@(18, 20330051920355, "URBANO", "GARCIA", "EVELYN", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "2APM", 2)
    ,@(19, 20330051920382, "VASQUEZ", "HERNANDEZ", "VICTOR MANUEL", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "2APM", 2)
    ,@(20, 20330051920357, "XILCAHUA", "TLAXCALA", "LUIS ANGEL", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "2APM", 2)
    ,@(21, 20330051920180, "RIVERA", "AGUILAR", "JAROMI YAJAIRA", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "2ARHM", 2)
    ,@(22, 20330051920039, "ALVAREZ", "CONCHE", "GUSTAVO", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "2BEM", 2)
    ,@(23, 20330051920040, "AMADOR", "PORRAS", "FRANCISCO ALAN", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "2BEM", 2)
    ,@(24, 20330051920097, "DE JESUS", "EVARISTO", "ALDAIR ALAN", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "2BEM", 2)
    ,@(25, 20330051920107, "SANCHEZ", "SANCHEZ", "ROSA ISELA", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "2BEM", 2)
    ,@(26, 20330051920286, "ALFONSO", "OSORIO", "AMERICA MICHELLE", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "2BLCM", 2)
    ,@(27, 20330051920295, "HERRERA", "CERON", "YAMILE", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "2BLCM", 2)
    ,@(28, 20330051920299, "LOPEZ", "MONTERROSAS", "MARIA MAGDALENA", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "2BLCM", 2)
    ,@(29, 20330051920301, "MARTINEZ", "HERNANDEZ", "VANESSA", "LECTURA, EXPRESIÓN ORAL Y ESCRITA II", "2BLCM", 2)
)

foreach ($row in $rescatables) {
    $r = $row[0]
    $ws4.Cells.Item($r, 2).Value = $row[2]
}
foreach ($row in $rescatables) {
    $r = $row[0]
    $ws4.Cells.Item($r, 3).Value = $row[3]
}
foreach ($row in $rescatables) {
    $r = $row[0]
    $ws4.Cells.Item($r, 4).Value = $row[4]
}
foreach ($row in $rescatables) {
    $r = $row[0]
    $ws4.Cells.Item($r, 1).Value = $row[1]
    $ws4.Cells.Item($r, 5).Value = $row[5]
    $ws4.Cells.Item($r, 6).Value = $row[6]
    $ws4.Cells.Item($r, 7).Value = $row[7]
}

Write-Host "Update complete"
